# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker/period detail table (rows 16-24, columns C:G) is re-sorted by
# worker, and two "Valor Mora" amounts are corrected:
#   - ELIANA ANDREA LORDUY TOBIO (periods 2202 / 2203): 1100000 -> 1623500
#   - JULIO ANDRES APRAEZ ESPAÑA (period 2210):          3600000 -> 3942000

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 16; Doc = "73182225";   Name = "POLICARPO MARTINEZ GARCIA";    Period = "2111"; Mora = 6057;  Valor = 1480000 },
    @{ Row = 17; Doc = "77094911";   Name = "JUSUE JOSE ONATE ROQUE";       Period = "2111"; Mora = 6057;  Valor = 1400000 },
    @{ Row = 18; Doc = "73168115";   Name = "LUIS ENRIQUE ABELLO MENDOZA";  Period = "2112"; Mora = 45760; Valor = 1000000 },
    @{ Row = 19; Doc = "73182225";   Name = "POLICARPO MARTINEZ GARCIA";    Period = "2112"; Mora = 36341; Valor = 1480000 },
    @{ Row = 20; Doc = "1047402132"; Name = "ELIANA ANDREA LORDUY TOBIO";   Period = "2202"; Mora = 44000; Valor = 1623500 },
    @{ Row = 21; Doc = "1047402132"; Name = "ELIANA ANDREA LORDUY TOBIO";   Period = "2203"; Mora = 44000; Valor = 1623500 },
    @{ Row = 22; Doc = "79498025";   Name = "JULIO ANDRES APRAEZ ESPAÑA";   Period = "2210"; Mora = 82666; Valor = 3942000 },
    @{ Row = 23; Doc = "9149836";    Name = "NELSON YAIR ARIZA ANZOATEGUI"; Period = "2412"; Mora = 80000; Valor = 2000000 },
    @{ Row = 24; Doc = "9149836";    Name = "NELSON YAIR ARIZA ANZOATEGUI"; Period = "2504"; Mora = 80000; Valor = 2000000 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("C$r").Value = $item.Doc
    $ws.Range("D$r").Value = $item.Name
    $ws.Range("E$r").Value = $item.Period
    $ws.Range("F$r").Value = $item.Mora
    $ws.Range("G$r").Value = $item.Valor
}
